$d = $word.ActiveDocument

# --- Create the three new character styles (wdStyleTypeCharacter = 2) ---
$ganStyle = $d.Styles.Add("GaNStyle", 2)
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.Bold = $true
$ganLinks.Font.Color = 8388608
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Underline = 1

# --- Apply GaNStyle to every "2022 Ημερομηνίες ..." run (4 occurrences) ---
$rng = $d.Content
$rng.Start = 0
$find = $rng.Find
$find.ClearFormatting()
$find.Text = "2022 Ημερομηνίες παρατήρησης για τον  Αστερισμός Πήγασος: 8-17 Οκτωβρίου, 7-16 Νοεμβρίου,"
$guard = 0
while ($find.Execute() -and $guard -lt 50) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $guard = $guard + 1
}

# --- Apply GaNParagraph to the "Συμμετέχετε ..." run ---
$rng2 = $d.Content
$rng2.Start = 0
$find2 = $rng2.Find
$find2.ClearFormatting()
$find2.Text = "Συμμετέχετε σε μία παγκόσμια καμπάνια για να παρατηρήσετε και να καταγράψετε τη φωτεινότητα των πιο αμυδρά ορατών άστρων σαν μέσο για την μέτρηση της Φωτορρύπανσης σε μία δεδομένη περιοχή. Με τον εντοπισμό και την παρατήρηση του  Αστερισμός Πήγασος στον νυχτερινό ουρανό καθώς και με την σύγκριση των ανωτέρω με τα διαγράμματα για τα μεγέθη των άστρων,  άνθρωποι από όλον τον κόσμο θα μάθουν πώς τα φώτα στην κοινότητά τους συμβάλλουν στην Φωτορρύπανση. Με την κατάθεση των πορισμάτων τους στην ιστοσελίδα θα δημιουργηθεί ένα αρχείο σχετικά με το τι μπορεί να δει κανείς στον νυχτερινό ουρανό."
$guard2 = 0
while ($find2.Execute() -and $guard2 -lt 50) {
    $rng2.Style = "GaNParagraph"
    $rng2.Collapse(0)
    $guard2 = $guard2 + 1
}

# --- Apply GaNLinks to the "Jan Hollan" credit run ---
$rng3 = $d.Content
$rng3.Start = 0
$find3 = $rng3.Find
$find3.ClearFormatting()
$find3.Text = "Τα διαγράμματα αυτού του αρχείου επιμελήθηκε ο Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$guard3 = 0
while ($find3.Execute() -and $guard3 -lt 50) {
    $rng3.Style = "GaNLinks"
    $rng3.Collapse(0)
    $guard3 = $guard3 + 1
}

Write-Host "Done"
